# Daily attendance processing - 2025-12-20 07:02:49
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Widen the "Status" column (column I, the 9th column) from 10 to 14 characters.
#    The engine's ColumnWidth setter adds a fixed padding of 5/6 (0.8333...) relative
#    to the raw stored width, so compensate to land exactly on 14.
$ws.Columns.Item(9).ColumnWidth = (14 - 5/6)

# 2. The "Recorded By" values that used to list the teacher before "System" now list
#    "System" first. Apply this sheet-wide since every occurrence gets reordered the
#    same way.
$ws.UsedRange.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")

# 3. Update the overview statistics: one more session (per group) is now officially
#    "missing" instead of "pending".
$ws.Range("L7").Value = 12
$ws.Range("L8").Value = 168

# 4. Each of the 12 group-statistics rows (one per B1 group) now shows one additional
#    missing session and one fewer pending session.
for ($r = 15; $r -le 26; $r++) {
    $ws.Range("P$r").Value = 1
    $ws.Range("Q$r").Value = 14
}

# 5. The "SURGERY SEMINAR/SLIDE" session dated 20/12/2025 has now passed for every
#    group, and since it was never recorded it flips from "Pending" (yellow) to
#    "Not Recorded" (red/pink), matching the legend's "Not Recorded" = Red status.
$notRecordedRows = @(17, 38, 59, 80, 100, 120, 140, 160, 180, 201, 222, 243)
$notRecordedColor = 12695295   # RGB(255, 182, 193) == the workbook's "Red" status fill

foreach ($r in $notRecordedRows) {
    $rowRange = $ws.Range("A" + $r + ":H" + $r)
    $rowRange.Interior.Color = $notRecordedColor

    $statusCell = $ws.Range("I$r")
    $statusCell.Value = "Not Recorded"
    $statusCell.Interior.Color = $notRecordedColor
}
